$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the _GoBack bookmark from the end of the "References"
#    paragraph to the end of the Author paragraph ("Joe Legner, P. E.,
#    S. E."). Bookmarks.Add() has a boundary bug in this runtime when
#    given a collapsed range that sits immediately before a paragraph
#    mark, so instead we delete the existing hidden bookmark outright
#    and splice a fresh bookmarkStart/bookmarkEnd pair straight into
#    the Author paragraph's OOXML via Range.InsertXML (preserving the
#    paragraph's existing identity attributes).
# ---------------------------------------------------------------------

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$authorPara = $d.Paragraphs.Item(2)

$authorXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="7BD58C44" w14:textId="77777777" w:rsidR="001D5706" w:rsidRDefault="00D207A2">
<w:pPr><w:pStyle w:val="Author"/></w:pPr>
<w:r><w:t>Joe Legner, P. E., S. E.</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$authorPara.Range.InsertXML($authorXml) | Out-Null

# ---------------------------------------------------------------------
# 2) Fix the Hyperlink character style's text color: it should follow
#    "auto" instead of the fixed theme accent color.
# ---------------------------------------------------------------------

$hyperlinkStyle = $d.Styles.Item("Hyperlink")
$hyperlinkStyle.Font.Color = -16777216   # wdColorAutomatic
